# Generate Report for Handoff
# Adds a new handed-off file (a97c9924-dc0f-48ef-bd6e-6cacb517eb39) as the
# next row on the Overview / zh-cn / de-de sheets, mirroring the existing
# 9a7cd8ad-... row that is already there.

$wb = $excel.ActiveWorkbook

$newGuid   = "a97c9924-dc0f-48ef-bd6e-6cacb517eb39"
$mdName    = "$newGuid.md"
$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/d29181dcc378652fbf83055f69a7de986dd29221/e2e/$mdName"

# ---------------------------------------------------------------------
# Overview sheet (row 3)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A3").Value = $mdName
$overview.Hyperlinks.Add($overview.Range("A3"), $mdUrl, "", "", $mdName) | Out-Null
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-29-18 14:29:58"

# ---------------------------------------------------------------------
# zh-cn sheet (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhXlfName = "$newGuid.4463417102ea53953bebea6f1432c9e0b33ed4ed.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7633c64d424c920535fef2cacb8521189e1a78f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"

$zhcn.Range("A3").Value = $mdName
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $mdUrl, "", "", $mdName) | Out-Null
$zhcn.Range("B3").Value = ".md"
$zhcn.Hyperlinks.Add($zhcn.Range("B3"), $mdUrl, "", "", ".md") | Out-Null
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = $zhXlfName
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), $zhXlfUrl, "", "", $zhXlfName) | Out-Null
$zhcn.Range("E3").Value = "2016-03-18 14:29:55"
$zhcn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("H3").Value = "0001-01-01 00:00:00"
$zhcn.Range("I3").Value = "Include"

# ---------------------------------------------------------------------
# de-de sheet (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$deXlfName = "$newGuid.4463417102ea53953bebea6f1432c9e0b33ed4ed.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd3516129b983ce73121dbf58942f63e05060d4f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$dede.Range("A3").Value = $mdName
$dede.Hyperlinks.Add($dede.Range("A3"), $mdUrl, "", "", $mdName) | Out-Null
$dede.Range("B3").Value = ".md"
$dede.Hyperlinks.Add($dede.Range("B3"), $mdUrl, "", "", ".md") | Out-Null
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = $deXlfName
$dede.Hyperlinks.Add($dede.Range("D3"), $deXlfUrl, "", "", $deXlfName) | Out-Null
$dede.Range("E3").Value = "2016-03-18 14:29:58"
$dede.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("H3").Value = "0001-01-01 00:00:00"
$dede.Range("I3").Value = "Include"
